$wb = $excel.ActiveWorkbook

# --- Sheet "Neg_Change": replace data rows 2-14 and append new rows 15-17 ---
$ws1 = $wb.Worksheets.Item("Neg_Change")

$ws1.Cells.Item(2,1).Value = "LT"
$ws1.Cells.Item(2,2).Value = 4001.1
$ws1.Cells.Item(2,3).Value = 4045.9
$ws1.Cells.Item(2,4).Value = 3980.2
$ws1.Cells.Item(2,5).Value = 4028
$ws1.Cells.Item(2,6).Value = 2328183
$ws1.Cells.Item(2,7).Value = 4959400
$ws1.Cells.Item(2,8).Value = -0.5305514780013711
$ws1.Cells.Item(2,9).Value = "LT"

$ws1.Cells.Item(3,1).Value = "DRREDDY"
$ws1.Cells.Item(3,2).Value = 1202.5
$ws1.Cells.Item(3,3).Value = 1205.5
$ws1.Cells.Item(3,4).Value = 1191
$ws1.Cells.Item(3,5).Value = 1198
$ws1.Cells.Item(3,6).Value = 4250157
$ws1.Cells.Item(3,7).Value = 8883497
$ws1.Cells.Item(3,8).Value = -0.5215671260991026
$ws1.Cells.Item(3,9).Value = "DRREDDY"

$ws1.Cells.Item(4,1).Value = "CGPOWER"
$ws1.Cells.Item(4,2).Value = 735
$ws1.Cells.Item(4,3).Value = 743.15
$ws1.Cells.Item(4,4).Value = 726.6
$ws1.Cells.Item(4,5).Value = 735.85
$ws1.Cells.Item(4,6).Value = 2560022
$ws1.Cells.Item(4,7).Value = 5090803
$ws1.Cells.Item(4,8).Value = -0.4971280562221717
$ws1.Cells.Item(4,9).Value = "CGPOWER"

$ws1.Cells.Item(5,1).Value = "ADANIGREEN"
$ws1.Cells.Item(5,2).Value = 1145
$ws1.Cells.Item(5,3).Value = 1159.5
$ws1.Cells.Item(5,4).Value = 1125.1
$ws1.Cells.Item(5,5).Value = 1135.9
$ws1.Cells.Item(5,6).Value = 5725034
$ws1.Cells.Item(5,7).Value = 11727335
$ws1.Cells.Item(5,8).Value = -0.5118213984677679
$ws1.Cells.Item(5,9).Value = "ADANIGREEN"

$ws1.Cells.Item(6,1).Value = "BPCL"
$ws1.Cells.Item(6,2).Value = 358
$ws1.Cells.Item(6,3).Value = 360.55
$ws1.Cells.Item(6,4).Value = 355.1
$ws1.Cells.Item(6,5).Value = 356.2
$ws1.Cells.Item(6,6).Value = 7532962
$ws1.Cells.Item(6,7).Value = 14807100
$ws1.Cells.Item(6,8).Value = -0.491260138717237
$ws1.Cells.Item(6,9).Value = "BPCL"

$ws1.Cells.Item(7,1).Value = "GLENMARK"
$ws1.Cells.Item(7,2).Value = 1883
$ws1.Cells.Item(7,3).Value = 1904.8
$ws1.Cells.Item(7,4).Value = 1876
$ws1.Cells.Item(7,5).Value = 1891.9
$ws1.Cells.Item(7,6).Value = 398054
$ws1.Cells.Item(7,7).Value = 854083
$ws1.Cells.Item(7,8).Value = -0.5339399098214108
$ws1.Cells.Item(7,9).Value = "GLENMARK"

$ws1.Cells.Item(8,1).Value = "EXIDEIND"
$ws1.Cells.Item(8,2).Value = 381.3
$ws1.Cells.Item(8,3).Value = 386.8
$ws1.Cells.Item(8,4).Value = 380
$ws1.Cells.Item(8,5).Value = 382.2
$ws1.Cells.Item(8,6).Value = 1778451
$ws1.Cells.Item(8,7).Value = 3586392
$ws1.Cells.Item(8,8).Value = -0.5041113743282943
$ws1.Cells.Item(8,9).Value = "EXIDEIND"

$ws1.Cells.Item(9,1).Value = "AUBANK"
$ws1.Cells.Item(9,2).Value = 878.25
$ws1.Cells.Item(9,3).Value = 882.75
$ws1.Cells.Item(9,4).Value = 872.7
$ws1.Cells.Item(9,5).Value = 876
$ws1.Cells.Item(9,6).Value = 710153
$ws1.Cells.Item(9,7).Value = 1730444
$ws1.Cells.Item(9,8).Value = -0.5896122613618239
$ws1.Cells.Item(9,9).Value = "AUBANK"

$ws1.Cells.Item(10,1).Value = "KPITTECH"
$ws1.Cells.Item(10,2).Value = 1170
$ws1.Cells.Item(10,3).Value = 1181.1
$ws1.Cells.Item(10,4).Value = 1161.4
$ws1.Cells.Item(10,5).Value = 1168
$ws1.Cells.Item(10,6).Value = 385752
$ws1.Cells.Item(10,7).Value = 879409
$ws1.Cells.Item(10,8).Value = -0.56135086177194
$ws1.Cells.Item(10,9).Value = "KPITTECH"

$ws1.Cells.Item(11,1).Value = "GMRAIRPORT"
$ws1.Cells.Item(11,2).Value = 95.2
$ws1.Cells.Item(11,3).Value = 95.94
$ws1.Cells.Item(11,4).Value = 93.7
$ws1.Cells.Item(11,5).Value = 93.83
$ws1.Cells.Item(11,6).Value = 12371529
$ws1.Cells.Item(11,7).Value = 26516655
$ws1.Cells.Item(11,8).Value = -0.5334430756820572
$ws1.Cells.Item(11,9).Value = "GMRAIRPORT"

$ws1.Cells.Item(12,1).Value = "INDUSTOWER"
$ws1.Cells.Item(12,2).Value = 367.75
$ws1.Cells.Item(12,3).Value = 371.35
$ws1.Cells.Item(12,4).Value = 362.8
$ws1.Cells.Item(12,5).Value = 364.05
$ws1.Cells.Item(12,6).Value = 6866183
$ws1.Cells.Item(12,7).Value = 15575932
$ws1.Cells.Item(12,8).Value = -0.5591799578991485
$ws1.Cells.Item(12,9).Value = "INDUSTOWER"

$ws1.Cells.Item(13,1).Value = "BLUESTARCO"
$ws1.Cells.Item(13,2).Value = 1965
$ws1.Cells.Item(13,3).Value = 1988
$ws1.Cells.Item(13,4).Value = 1932
$ws1.Cells.Item(13,5).Value = 1932.5
$ws1.Cells.Item(13,6).Value = 127593
$ws1.Cells.Item(13,7).Value = 256389
$ws1.Cells.Item(13,8).Value = -0.502346044487088
$ws1.Cells.Item(13,9).Value = "BLUESTARCO"

$ws1.Cells.Item(14,1).Value = "AMBER"
$ws1.Cells.Item(14,2).Value = 8142.5
$ws1.Cells.Item(14,3).Value = 8176.5
$ws1.Cells.Item(14,4).Value = 8015
$ws1.Cells.Item(14,5).Value = 8033.5
$ws1.Cells.Item(14,6).Value = 164814
$ws1.Cells.Item(14,7).Value = 362118
$ws1.Cells.Item(14,8).Value = -0.5448610673868739
$ws1.Cells.Item(14,9).Value = "AMBER"

$ws1.Cells.Item(15,1).Value = "CAMS"
$ws1.Cells.Item(15,2).Value = 3965
$ws1.Cells.Item(15,3).Value = 3975
$ws1.Cells.Item(15,4).Value = 3914.4
$ws1.Cells.Item(15,5).Value = 3926.8
$ws1.Cells.Item(15,6).Value = 266911
$ws1.Cells.Item(15,7).Value = 663840
$ws1.Cells.Item(15,8).Value = -0.5979287177633165
$ws1.Cells.Item(15,9).Value = "CAMS"

$ws1.Cells.Item(16,1).Value = "NBCC"
$ws1.Cells.Item(16,2).Value = 118.3
$ws1.Cells.Item(16,3).Value = 119.81
$ws1.Cells.Item(16,4).Value = 117.05
$ws1.Cells.Item(16,5).Value = 117.2
$ws1.Cells.Item(16,6).Value = 11263219
$ws1.Cells.Item(16,7).Value = 24956240
$ws1.Cells.Item(16,8).Value = -0.5486812516629107
$ws1.Cells.Item(16,9).Value = "NBCC"

$ws1.Cells.Item(17,1).Value = "KFINTECH"
$ws1.Cells.Item(17,2).Value = 1102
$ws1.Cells.Item(17,3).Value = 1118
$ws1.Cells.Item(17,4).Value = 1091.2
$ws1.Cells.Item(17,5).Value = 1097.8
$ws1.Cells.Item(17,6).Value = 652110
$ws1.Cells.Item(17,7).Value = 1359143
$ws1.Cells.Item(17,8).Value = -0.5202050115403604
$ws1.Cells.Item(17,9).Value = "KFINTECH"

# --- Sheet "Pos_Change": replace data rows 2-10 and append new row 11 ---
$ws2 = $wb.Worksheets.Item("Pos_Change")

$ws2.Cells.Item(2,1).Value = "ITC"
$ws2.Cells.Item(2,2).Value = 420.5
$ws2.Cells.Item(2,3).Value = 426.4
$ws2.Cells.Item(2,4).Value = 418.2
$ws2.Cells.Item(2,5).Value = 419.95
$ws2.Cells.Item(2,6).Value = 18599864
$ws2.Cells.Item(2,7).Value = 12938674
$ws2.Cells.Item(2,8).Value = 0.4375401992507115
$ws2.Cells.Item(2,9).Value = "ITC"

$ws2.Cells.Item(3,1).Value = "TRENT"
$ws2.Cells.Item(3,2).Value = 4745
$ws2.Cells.Item(3,3).Value = 4748.9
$ws2.Cells.Item(3,4).Value = 4685.2
$ws2.Cells.Item(3,5).Value = 4697.3
$ws2.Cells.Item(3,6).Value = 498496
$ws2.Cells.Item(3,7).Value = 346569
$ws2.Cells.Item(3,8).Value = 0.4383744651137292
$ws2.Cells.Item(3,9).Value = "TRENT"

$ws2.Cells.Item(4,1).Value = "KOTAKBANK"
$ws2.Cells.Item(4,2).Value = 2146.9
$ws2.Cells.Item(4,3).Value = 2146.9
$ws2.Cells.Item(4,4).Value = 2100
$ws2.Cells.Item(4,5).Value = 2105.5
$ws2.Cells.Item(4,6).Value = 3786254
$ws2.Cells.Item(4,7).Value = 2487693
$ws2.Cells.Item(4,8).Value = 0.5219940724197077
$ws2.Cells.Item(4,9).Value = "KOTAKBANK"

$ws2.Cells.Item(5,1).Value = "OBEROIRLTY"
$ws2.Cells.Item(5,2).Value = 1760
$ws2.Cells.Item(5,3).Value = 1794.9
$ws2.Cells.Item(5,4).Value = 1759
$ws2.Cells.Item(5,5).Value = 1777.1
$ws2.Cells.Item(5,6).Value = 1586130
$ws2.Cells.Item(5,7).Value = 1074751
$ws2.Cells.Item(5,8).Value = 0.4758116065953881
$ws2.Cells.Item(5,9).Value = "OBEROIRLTY"

$ws2.Cells.Item(6,1).Value = "TORNTPOWER"
$ws2.Cells.Item(6,2).Value = 1316
$ws2.Cells.Item(6,3).Value = 1323
$ws2.Cells.Item(6,4).Value = 1301.5
$ws2.Cells.Item(6,5).Value = 1317.6
$ws2.Cells.Item(6,6).Value = 269394
$ws2.Cells.Item(6,7).Value = 172378
$ws2.Cells.Item(6,8).Value = 0.5628096392811147
$ws2.Cells.Item(6,9).Value = "TORNTPOWER"

$ws2.Cells.Item(7,1).Value = "DALBHARAT"
$ws2.Cells.Item(7,2).Value = 2100.6
$ws2.Cells.Item(7,3).Value = 2109.9
$ws2.Cells.Item(7,4).Value = 2087.4
$ws2.Cells.Item(7,5).Value = 2100
$ws2.Cells.Item(7,6).Value = 304566
$ws2.Cells.Item(7,7).Value = 213064
$ws2.Cells.Item(7,8).Value = 0.4294578154920587
$ws2.Cells.Item(7,9).Value = "DALBHARAT"

$ws2.Cells.Item(8,1).Value = "CUMMINSIND"
$ws2.Cells.Item(8,2).Value = 4392
$ws2.Cells.Item(8,3).Value = 4399.7
$ws2.Cells.Item(8,4).Value = 4333
$ws2.Cells.Item(8,5).Value = 4351
$ws2.Cells.Item(8,6).Value = 677463
$ws2.Cells.Item(8,7).Value = 444822
$ws2.Cells.Item(8,8).Value = 0.522997963230236
$ws2.Cells.Item(8,9).Value = "CUMMINSIND"

$ws2.Cells.Item(9,1).Value = "OFSS"
$ws2.Cells.Item(9,2).Value = 8633.5
$ws2.Cells.Item(9,3).Value = 8637
$ws2.Cells.Item(9,4).Value = 8500
$ws2.Cells.Item(9,5).Value = 8525
$ws2.Cells.Item(9,6).Value = 111369
$ws2.Cells.Item(9,7).Value = 75010
$ws2.Cells.Item(9,8).Value = 0.4847220370617251
$ws2.Cells.Item(9,9).Value = "OFSS"

$ws2.Cells.Item(10,1).Value = "INOXWIND"
$ws2.Cells.Item(10,2).Value = 154.55
$ws2.Cells.Item(10,3).Value = 157.26
$ws2.Cells.Item(10,4).Value = 154.26
$ws2.Cells.Item(10,5).Value = 154.77
$ws2.Cells.Item(10,6).Value = 6011244
$ws2.Cells.Item(10,7).Value = 4118297
$ws2.Cells.Item(10,8).Value = 0.4596431486121569
$ws2.Cells.Item(10,9).Value = "INOXWIND"

$ws2.Cells.Item(11,1).Value = "NUVAMA"
$ws2.Cells.Item(11,2).Value = 7247
$ws2.Cells.Item(11,3).Value = 7399
$ws2.Cells.Item(11,4).Value = 6999
$ws2.Cells.Item(11,5).Value = 7115
$ws2.Cells.Item(11,6).Value = 142339
$ws2.Cells.Item(11,7).Value = 101601
$ws2.Cells.Item(11,8).Value = 0.400960620466334
$ws2.Cells.Item(11,9).Value = "NUVAMA"

